$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 2.647218666666667
$ws.Range("H2").Value = 7.941656
$ws.Range("I2").Value = 0.4640059894538357
$ws.Range("J2").Value = 0.4640059894538356
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 1.195417
$ws.Range("N2").Value = 3.586251
$ws.Range("O2").Value = 0.724690891256891
$ws.Range("P2").Value = 0.7246908912568911
$ws.Range("Q2").Value = 3.164530196850667
$ws.Range("R2").Value = 28.480771771656
$ws.Range("S2").Value = 0.3362609140458357
$ws.Range("T2").Value = 0.3362609140458358

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 2.647218666666667
$ws.Range("H3").Value = 7.941656
$ws.Range("I3").Value = 0.4640059894538357
$ws.Range("J3").Value = 0.4640059894538356
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.4541373333333333
$ws.Range("N3").Value = 1.362412
$ws.Range("O3").Value = 0.275309108743109
$ws.Range("P3").Value = 0.275309108743109
$ws.Range("Q3").Value = 1.202200826030222
$ws.Range("R3").Value = 10.819807434272
$ws.Range("S3").Value = 0.1277450754079999
$ws.Range("T3").Value = 0.1277450754079999

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 3.057920333333333
$ws.Range("H4").Value = 9.173760999999999
$ws.Range("I4").Value = 0.5359940105461642
$ws.Range("J4").Value = 0.5359940105461642
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 1.195417
$ws.Range("N4").Value = 3.586251
$ws.Range("O4").Value = 0.724690891256891
$ws.Range("P4").Value = 0.7246908912568911
$ws.Range("Q4").Value = 3.655489951112334
$ws.Range("R4").Value = 32.899409560011
$ws.Range("S4").Value = 0.3884299772110552
$ws.Range("T4").Value = 0.3884299772110553

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 3.057920333333333
$ws.Range("H5").Value = 9.173760999999999
$ws.Range("I5").Value = 0.5359940105461642
$ws.Range("J5").Value = 0.5359940105461642
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.4541373333333333
$ws.Range("N5").Value = 1.362412
$ws.Range("O5").Value = 0.275309108743109
$ws.Range("P5").Value = 0.275309108743109
$ws.Range("Q5").Value = 1.388715785725778
$ws.Range("R5").Value = 12.498442071532
$ws.Range("S5").Value = 0.147564033335109
$ws.Range("T5").Value = 0.147564033335109
